$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.853.67"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.765.68"
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'327.28"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4468"
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("D8").Value = "'0.3541"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "'0.07432"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "'42.09"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("D13").Value = "'20.79"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'6.013"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "'7.186"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "1.769.64"
$ws.Range("D17").Value = "'92.83"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'0.06435"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'5.775"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "27.891.89"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'11.27"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.129"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "'162.75"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'20.15"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "1.973.93"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "'2.160"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("D30").Value = "'124.79"
$ws.Range("D31").Value = "'1.096"
$ws.Range("E31").Value = "  +3.98%  "
$ws.Range("D32").Value = "'0.09132"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "'3.650"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'5.547"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'11.81"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "'0.02288"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "'0.06091"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'0.2087"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6292"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.953"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'1.179"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "'1.389"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "'7.919"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "'3.732"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "'0.5838"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'121.93"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "'1.945"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "'0.06904"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'1.134"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'72.61"
$ws.Range("E51").Value = "  +0.52%  "
